$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: D-column price values that are valid numeric literals (e.g. "1.000",
# "0.9999") are written with a leading apostrophe so Excel stores them as
# literal text (matching the source General-formatted text cells) instead of
# silently normalizing them into numbers (which would turn "1.000" into 1).

# Row 2
$ws.Range("D2").Value = "29.329.33"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").Value = "1.861.83"

# Row 4
$ws.Range("D4").Value = "'0.9999"

# Row 5
$ws.Range("D5").Value = "'0.7056"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("D6").Value = "'238.38"
$ws.Range("E6").Value = "  +0.29%  "

# Row 8
$ws.Range("D8").Value = "'0.07899"
$ws.Range("E8").Value = "  +2.11%  "

# Row 9
$ws.Range("D9").Value = "'0.3055"
$ws.Range("E9").Value = "  +0.19%  "

# Row 10
$ws.Range("D10").Value = "'24.84"
$ws.Range("E10").Value = "  +6.78%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.08192"
$ws.Range("E11").Value = "  -0.01%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.964.62"
$ws.Range("E12").Value = "  +5.47%  "

# Row 13
$ws.Range("D13").Value = "'5.231"
$ws.Range("E13").Value = "  +1.06%  "

# Row 14
$ws.Range("D14").Value = "'0.7177"
$ws.Range("E14").Value = "  -0.08%  "

# Row 15
$ws.Range("D15").Value = "'89.65"
$ws.Range("E15").Value = "  +0.52%  "

# Row 16
$ws.Range("D16").Value = "29.706.92"
$ws.Range("E16").Value = "  +1.62%  "

# Row 17
$ws.Range("D17").Value = "'5.857"
$ws.Range("E17").Value = "  +1.34%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007827"
$ws.Range("E18").Value = "  +1.20%  "

# Row 19
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.249.31"
$ws.Range("E19").Value = "  +7.06%  "

# Row 20
$ws.Range("E20").Value = "  -0.74%  "

# Row 21
$ws.Range("D21").Value = "'238.95"
$ws.Range("E21").Value = "  +0.72%  "

# Row 22
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("D24").Value = "'7.607"
$ws.Range("E24").Value = "  +2.23%  "

# Row 25
$ws.Range("D25").Value = "'163.02"
$ws.Range("E25").Value = "  +0.66%  "

# Row 26
$ws.Range("E26").Value = "  -1.01%  "

# Row 27
$ws.Range("E27").Value = "  -2.60%  "

# Row 28
$ws.Range("D28").Value = "'18.12"
$ws.Range("E28").Value = "  +0.49%  "

# Row 29
$ws.Range("D29").Value = "'1.905"
$ws.Range("E29").Value = "  -5.92%  "

# Row 30
$ws.Range("D30").Value = "'1.375"
$ws.Range("E30").Value = "  -3.96%  "

# Row 31
$ws.Range("D31").Value = "'1.479"
$ws.Range("E31").Value = "  -0.60%  "

# Row 32
$ws.Range("D32").Value = "'4.326"
$ws.Range("E32").Value = "  -2.33%  "

# Row 33
$ws.Range("D33").Value = "'4.063"
$ws.Range("E33").Value = "  +0.18%  "

# Row 34
$ws.Range("D34").Value = "'0.05187"
$ws.Range("E34").Value = "  -0.84%  "

# Row 35
$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = "  +0.67%  "

# Row 36
$ws.Range("D36").Value = "'0.7106"
$ws.Range("E36").Value = "  +0.35%  "

# Row 37
$ws.Range("D37").Value = "'1.007"
$ws.Range("E37").Value = "  +0.73%  "

# Row 38
$ws.Range("D38").Value = "'2.675"
$ws.Range("E38").Value = "  +0.25%  "

# Row 39
$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = "  +0.49%  "

# Row 40
$ws.Range("D40").Value = "'2.689"
$ws.Range("E40").Value = "  -1.34%  "

# Row 41
$ws.Range("D41").Value = "1.179.78"
$ws.Range("E41").Value = "  +3.43%  "

# Row 42
$ws.Range("D42").Value = "'0.9220"
$ws.Range("E42").Value = "  -1.09%  "

# Row 43
$ws.Range("D43").Value = "'6.039"
$ws.Range("E43").Value = "  +2.11%  "

# Row 44
$ws.Range("D44").Value = "'72.09"
$ws.Range("E44").Value = "  +1.68%  "

# Row 45
$ws.Range("D45").Value = "'0.4276"
$ws.Range("E45").Value = "  -0.08%  "

# Row 46
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.11%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'102.25"
$ws.Range("E47").Value = "  -1.32%  "

# Row 48
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.133.14"
$ws.Range("E48").Value = "  +6.48%  "

# Row 49
$ws.Range("D49").Value = "'0.5347"
$ws.Range("E49").Value = "  -1.86%  "

# Row 50
$ws.Range("D50").Value = "'1.761"
$ws.Range("E50").Value = "  -1.83%  "

# Row 51
$ws.Range("D51").Value = "'9.192"
$ws.Range("E51").Value = "  +0.03%  "
